$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Auto")

# Week 9 (column M) results for each team row (3-34).
# Row 11 (Jeff / Chargers) also gets a correction to its week 8 (column L)
# cumulative value: 4 -> 3, since week 9 brings it to 4.
$weekNine = @{
    3  = 1
    4  = 2
    5  = 3
    6  = 4
    7  = 6
    8  = 0
    9  = 3
    10 = 3
    11 = 4
    12 = 6
    13 = 4
    14 = 7
    15 = 4
    16 = 4
    17 = 6
    18 = 5
    19 = 2
    20 = 3
    21 = 5
    22 = 4
    23 = 3
    24 = 7
    25 = 7
    26 = 3
    27 = 4
    28 = 4
    29 = 4
    30 = 5
    31 = 4
    32 = 5
    33 = 4
    34 = 5
}

# Correct the pre-existing week 8 (L11) value for row 11.
$ws.Range("L11").Value = 3

foreach ($row in 3..34) {
    $ws.Cells.Item($row, 13).Value = $weekNine[$row]
}

# Extend the per-row CONCATENATE formula in column A to include the new
# week 9 column (M).
$ws.Range("A3").Formula = '=CONCATENATE("[''",$C3,"'', ","''",$D3,"'' ,",$E3,",",$F3,",",$G3,",",$H3,",",$I3,",",$J3,",",$K3,",",$L3,",",$M3,"],")'
$ws.Range("A4:A34").Formula = '=CONCATENATE("[''",$C4,"'', ","''",$D4,"'' ,",$E4,",",$F4,",",$G4,",",$H4,",",$I4,",",$J4,",",$K4,",",$L4,",",$M4,"],")'

# Extend the SUMIF totals at the bottom of the sheet (rows 36-39) to cover
# the new week 9 column, matching the right-aligned style already used by
# the rest of that summary block (columns F:L).
$ws.Range("M36:M39").HorizontalAlignment = -4152
$ws.Range("M36:M38").Formula = '=SUMIF($C$3:$C$34,$D36,M$3:M$34)'
$ws.Range("M39").Formula = '=SUMIF($C$3:$C$34,$D$39,M$3:M$34)'

# Update the frozen-pane scroll position and active selection to reflect
# the newly revealed week 9 column.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("L1").Select() | Out-Null
$win.FreezePanes = $true
$ws.Range("M2").Select() | Out-Null
